$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was M, now B)
$ws.Range("A2").Value = "B"
$ws.Range("B2").Value = 0.9714285714285714
$ws.Range("C2").Value = 0.9444444444444444
$ws.Range("D2").Value = 0.9577464788732394
$ws.Range("E2").Value = 36

# Row 3 (was B, now M)
$ws.Range("A3").Value = "M"
$ws.Range("B3").Value = 0.9090909090909091
$ws.Range("C3").Value = 0.9523809523809523
$ws.Range("D3").Value = 0.9302325581395349
$ws.Range("E3").Value = 21

# Row 4 (accuracy)
$ws.Range("B4").Value = 0.9473684210526315
$ws.Range("C4").Value = 0.9473684210526315
$ws.Range("D4").Value = 0.9473684210526315
$ws.Range("E4").Value = 0.9473684210526315

# Row 5 (macro avg)
$ws.Range("B5").Value = 0.9402597402597402
$ws.Range("C5").Value = 0.9484126984126984
$ws.Range("D5").Value = 0.9439895185063871
$ws.Range("E5").Value = 57

# Row 6 (weighted avg)
$ws.Range("B6").Value = 0.9484620642515378
$ws.Range("C6").Value = 0.9473684210526315
$ws.Range("D6").Value = 0.9476097712345061
$ws.Range("E6").Value = 57
